# Section 2.5 "Design and Implementation Constraints" is appended at the very
# end of the document body (right before the trailing sectPr), immediately
# after the last existing paragraph ("No hosting plan yet for future backend").
#
# We build the new content as a raw OOXML fragment (rather than typing text
# and then reformatting) so the inserted paragraphs get exactly the pPr/rPr
# combinations from the target revision, with no unwanted inheritance (e.g.
# bullet/numbering) from the preceding list paragraph.

$d = $word.ActiveDocument

$anchor = $d.Content
$found = $anchor.Find.Execute("No hosting plan yet for future backend", $false, $false,
                               $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph 'No hosting plan yet for future backend' not found"
}

$insertAt = $d.Range($anchor.End, $anchor.End)

$newParagraphsXml = (
    '<w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0"/>' +
    '<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:bCs/><w:sz w:val="22"/>' +
    '<w:szCs w:val="22"/></w:rPr></w:pPr></w:p>' +

    '<w:p><w:pPr><w:pStyle w:val="Heading2"/><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr>' +
    '<w:t>Design and Implementation Constraints</w:t></w:r></w:p>' +

    '<w:p><w:pPr><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr></w:pPr>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr><w:t>Syarti</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr>' +
    '<w:t xml:space="preserve"> project is being developed solely for educational purposes and is ' +
    'currently subject to several design and implementation constraints. While there are no restrictions ' +
    'on the use of backend technologies, the team lacks the required experience in backend development, ' +
    'databases, and API integration. As a result, the system is currently limited to a front-end-only ' +
    'prototype hosted via GitHub Pages by choice, with no external limitations.</w:t></w:r></w:p>' +

    '<w:p><w:pPr><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Times New Roman"/></w:rPr>' +
    '<w:t>The application is designed for Arabic-speaking users only, and no multilingual or accessibility ' +
    'support is currently implemented. Due to the academic nature of the project, there are no requirements ' +
    'to follow privacy laws, rental regulations, or data protection practices, and the system does not ' +
    'store or handle any user data.</w:t></w:r></w:p>' +

    '<w:p><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:before="0" w:beforeAutospacing="0"/>' +
    '<w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi"/><w:bCs/><w:sz w:val="22"/>' +
    '<w:szCs w:val="22"/></w:rPr></w:pPr></w:p>'
)

$flatOpc = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
$newParagraphsXml
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$insertAt.InsertXML($flatOpc)

Write-Host "Inserted 'Design and Implementation Constraints' section. Paragraph count now:" $d.Paragraphs.Count
